# Standard User name change in Contacts test data files - 6th Mar 2024
#
# The "Users" sheet's StdUser record (row 2) is renamed from
# "Drew Koecher" to "Ayati Arvind". The previously-active "Relationship"
# sheet's lingering selection is also moved (C21) and the "Users" sheet
# becomes the active tab/selection (E6) to reflect where the editor last
# worked before saving.

$wb = $excel.ActiveWorkbook

# Touch the Relationship sheet's selection first (it was the tab active
# before this edit), then move off of it so it's no longer the active tab.
$wsRelationship = $wb.Worksheets.Item("Relationship")
$wsRelationship.Activate()
$wsRelationship.Range("C21").Select()

# Make the actual data edit on the Users sheet.
$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Activate()
$wsUsers.Range("A2").Value = "Ayati Arvind"
$wsUsers.Range("E6").Select()
